# Add two new progress rows (陳家畯 / 12-16 and 12-23 entries) above the
# existing 羅致遠 row, pushing it from row 8 down to row 10, then restore
# the filter/defined-name bookkeeping that Excel re-wrote on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 8 (the 羅致遠 row), which
# shifts it down to row 10 and inherits row 7's formatting as a starting
# point; we overwrite the style explicitly for every cell below anyway.
$ws.Rows("8:9").Insert()

# --- Row 8: 陳家畯, filled-in 12/16 status -------------------------------
$ws.Range("A8").Value = 43815
$ws.Range("A8").NumberFormat = "mm-dd-yy"
$ws.Range("B8").Value = "陳家畯"
$ws.Range("C8").Value = "資料蒐集"
$ws.Range("D8").Value = "價量模擬"
$ws.Range("E8").Value = "參與12/15之討論，協助執行程式碼"
$ws.Range("F8").Value = "完成價量模型模擬"

# --- Row 9: 陳家畯, filled-in 12/23 status -------------------------------
$ws.Range("A9").Value = 43822
$ws.Range("A9").NumberFormat = "mm-dd-yy"
$ws.Range("B9").Value = "陳家畯"
$ws.Range("C9").Value = "資料蒐集"
$ws.Range("D9").Value = "價量模擬"
$ws.Range("E9").Value = "參與12/18之討論，協助執行程式碼結果"
$ws.Range("F9").Value = "完成期末影片剪輯"

# Match the column styles used throughout the sheet: column E wraps text
# (style 6), the rest align to the top without wrapping (style 5), and
# column A keeps the date number format (style 4).
$ws.Range("A8:A9").VerticalAlignment = -4160
$ws.Range("A8:A9").WrapText = $false

$ws.Range("B8:D9").VerticalAlignment = -4160
$ws.Range("B8:D9").WrapText = $false
$ws.Range("F8:F9").VerticalAlignment = -4160
$ws.Range("F8:F9").WrapText = $false

$ws.Range("E8:E9").VerticalAlignment = -4160
$ws.Range("E8:E9").WrapText = $true

# Row heights observed in the saved workbook.
$ws.Rows("8").RowHeight = 33.5
$ws.Rows("9").RowHeight = 31

# The used range grew; update the hidden _FilterDatabase defined name
# accordingly (Excel re-stamps this on save even though the visible
# AutoFilter range itself was left untouched at A1:F8).
$wb.Names("_xlnm._FilterDatabase").RefersTo = "=工作進度!`$A`$1:`$F`$12"
